# Quarterly financials update: a new quarter (period ending 2018-09-30)
# is inserted as a new column D, shifting the existing eight quarters
# (previously D:K) one column to the right (now E:L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before the current column D; this shifts the
# existing D:K data block to E:L and keeps per-column formatting/width.
$ws.Range("D1").EntireColumn.Insert()

# --- Income Statement ---------------------------------------------------
$ws.Range("D7").Value = 43373   # Period Ending (9/30/2018)
$ws.Range("D8").Value = 7900    # Total Revenue
$ws.Range("D9").Value = 3000    # Cost of Revenue
$ws.Range("D10").Value = 4900   # Gross Profit
# Row 11 "Operating Expenses" header has no value cells to fill.
$ws.Range("D12").Value = 200    # Research Development
$ws.Range("D13").Value = 0      # Selling General and Administrative
$ws.Range("D14").Value = 0      # Non Recurring
$ws.Range("D15").Value = 0      # Others
# Row 16 is a blank spacer row.
$ws.Range("D17").Value = 7900   # Total Operating Expenses
$ws.Range("D18").Value = 0      # Operating Income or Loss
# Row 19 "Income from Continuing Operations" header has no value cells.
$ws.Range("D20").Value = -100   # Total Other Income/Expenses Net
$ws.Range("D21").Value = 1100   # Earnings Before Interest And Taxes
$ws.Range("D22").Value = 200    # Interest Expense
$ws.Range("D23").Value = -400   # Income Before Tax
$ws.Range("D24").Value = -100   # Income Tax Expense
$ws.Range("D25").Value = 0      # Minority Interest
$ws.Range("D26").Value = -300   # Income After Tax
$ws.Range("D27").Value = -300   # Net Income From Continuing Ops
$ws.Range("D28").Value = 0      # Non-recurring Events
$ws.Range("D29").Value = "NA"   # Discontinued Operations
$ws.Range("D30").Value = 0      # Extraordinary Items
$ws.Range("D31").Value = 0      # Effect Of Accounting Changes
$ws.Range("D32").Value = 100    # Other Items
$ws.Range("D33").Value = -300   # Net Income
$ws.Range("D34").Value = 0      # Preferred Stock And Other Adjustments
$ws.Range("D35").Value = -300   # Net Income Applicable To Common Shares

# --- Balance Sheet --------------------------------------------------------
$ws.Range("D38").Value = 43373  # Period Ending (9/30/2018)
# Row 39 "Assets" / Row 40 "Current Assets" headers have no value cells.
$ws.Range("D41").Value = 15900  # Cash And Cash Equivalents
$ws.Range("D42").Value = 0      # Short Term Investments
$ws.Range("D43").Value = 2700   # Net Receivables
$ws.Range("D44").Value = 2500   # Inventory
$ws.Range("D45").Value = 700    # Other Current Assets
$ws.Range("D46").Value = 21800  # Total Current Assets
$ws.Range("D47").Value = 0      # Long Term Investments
$ws.Range("D48").Value = 5700   # Property Plant and Equipment
$ws.Range("D49").Value = 18700  # Goodwill
$ws.Range("D50").Value = 0      # Intangible Assets
$ws.Range("D51").Value = 0      # Accumulated Amortization
$ws.Range("D52").Value = 0      # Other Assets
$ws.Range("D53").Value = 0      # Deferred Long Term Asset Charges
$ws.Range("D54").Value = 46200  # Total Assets
# Row 55 "Liabilities" / Row 56 "Current Liabilities" headers have no value cells.
$ws.Range("D57").Value = 1700   # Accounts Payable
$ws.Range("D58").Value = 0      # Short/Current Long Term Debt
$ws.Range("D59").Value = 3100   # Other Current Liabilities
$ws.Range("D60").Value = 4800   # Total Current Liabilities
$ws.Range("D61").Value = 7400   # Long Term Debt
$ws.Range("D62").Value = 700    # Other Liabilities
$ws.Range("D63").Value = 0      # Deferred Long Term Liability Charges
$ws.Range("D64").Value = 0      # Minority Interest
$ws.Range("D65").Value = 0      # Negative Goodwill
$ws.Range("D66").Value = 12800  # Total Liabilities
# Row 67 "Stockholders' Equity" header has no value cells.
$ws.Range("D68").Value = 0      # Misc Stocks Options Warrants
$ws.Range("D69").Value = 0      # Redeemable Preferred Stock
$ws.Range("D70").Value = 0      # Preferred Stock
$ws.Range("D71").Value = 0      # Common Stock
$ws.Range("D72").Value = -233500 # Retained Earnings
$ws.Range("D73").Value = 0      # Treasury Stock
$ws.Range("D74").Value = 0      # Capital Surplus
$ws.Range("D75").Value = 0      # Other Stockholder Equity
$ws.Range("D76").Value = 33400  # Total Stockholder Equity
$ws.Range("D77").Value = 0      # Net Tangible Assets

# --- Cash Flow Statement ---------------------------------------------------
$ws.Range("D80").Value = 43373  # Period Ending (9/30/2018)
$ws.Range("D81").Value = -300   # Net Income
# Row 82 header has no value cells.
$ws.Range("D83").Value = 1300   # Depreciation
$ws.Range("D84").Value = 0      # Adjustments To Net Income
$ws.Range("D85").Value = 0      # Changes In Accounts Receivables
$ws.Range("D86").Value = 0      # Changes In Liabilities
$ws.Range("D87").Value = 0      # Changes In Inventories
$ws.Range("D88").Value = 0      # Changes In Other Operating Activities
$ws.Range("D89").Value = 1900   # Total Cash Flow From Operating Activities
# Row 90 header has no value cells.
$ws.Range("D91").Value = 0      # Capital Expenditures
$ws.Range("D92").Value = 0      # Investments
$ws.Range("D93").Value = 0      # Other Cashflows from Investing Activities
$ws.Range("D94").Value = -400   # Total Cash Flows From Investing Activities
# Row 95 header has no value cells.
$ws.Range("D96").Value = 0      # Dividends Paid
$ws.Range("D97").Value = 0      # Sale Purchase of Stock
$ws.Range("D98").Value = 0      # Net Borrowings
$ws.Range("D99").Value = 0      # Other Cash Flows from Financing Activities
$ws.Range("D100").Value = -100  # Total Cash Flows From Financing Activities
$ws.Range("D101").Value = "NA"  # Effect Of Exchange Rate Changes
$ws.Range("D102").Value = 1400  # Change In Cash and Cash Equivalents
